# "Chromatic Dispersion first implementation"
#
# 1) The deck's Date placeholder (an auto "datetimeFigureOut" field) is
#    re-cached from 11/10/2020 -> 11/11/2020 on the Slide Master and on
#    every Slide Layout (this happens automatically whenever PowerPoint
#    re-saves the file on a later day). We touch every "Date Placeholder*"
#    shape we can find, on the master and on each custom layout.
# 2) On slide 9 ("Upscaling and FFT"), the three bullet paragraphs in the
#    content placeholder were re-touched by the author while working on
#    the Chromatic Dispersion write-up (the runs are re-typed with the
#    same text, which is how PowerPoint clears the "needs reproofing"
#    state on a run).

$p = $ppt.ActivePresentation

# --- 1) Refresh the cached "today" date field text everywhere it lives ---

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide Master
Update-DatePlaceholder $p.SlideMaster.Shapes "11/11/2020"

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "11/11/2020"
}

# --- 2) Re-touch the three bullets on slide 9 ---

$slide9 = $p.Slides.Item(9)
$body = $slide9.Shapes.Item(3)
$tr = $body.TextFrame.TextRange

$paraTexts = @(
    "Upscaled the symbols by 8, plotted the spectrum before and after modulation",
    "Compared the spectrum with the spectrum of the symbols without upscaling",
    "Learning about the Filter Bank method for spectrum plotting"
)

for ($i = 1; $i -le $paraTexts.Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.Text = $paraTexts[$i - 1]
}
